$wb = $excel.ActiveWorkbook

# --- addListItem sheet: change the eli_text value in A2 from "Userseven"
#     to "Usereight" (the dependent C2 formula "=A2" recalculates itself). ---
$wsAddListItem = $wb.Worksheets.Item("addListItem")
$wsAddListItem.Range("A2").Value = "Usereight"

# --- createUser sheet: bump the increment value in A2 from 1019 to 1023
#     (the dependent CONCAT formulas in B2/F2 recalculate themselves). ---
$wsCreateUser = $wb.Worksheets.Item("createUser")
$wsCreateUser.Range("A2").Value = 1023

# --- Switch the active/selected tab to "addListItem" (it was previously
#     sitting on "setHpClinicDiary"). ---
$wsAddListItem.Activate()
